$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 27 ("09876543", blank birthday, 0 points) down to new row 28,
# preserving the original text value (with leading zero) via copy/paste,
# then set the new row's points to 120.
$ws.Range("A27:C27").Copy()
$ws.Range("A28").PasteSpecial()
$ws.Cells.Item(28, 2).Value = ""
$ws.Cells.Item(28, 3).Value = 120

# Row 27's phone number becomes a plain numeric value (leading zero dropped).
$ws.Cells.Item(27, 1).Value = 9876543
